# Update "想去人数" (number of people interested) counts for a few rows
# in the "展览" sheet and the "全部类型" sheet, as reflected in the
# regenerated output data (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9207
$ws1.Range("F3").Value = 205
$ws1.Range("F4").Value = 479
$ws1.Range("F5").Value = 457

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9207
$ws4.Range("F3").Value = 205
$ws4.Range("F4").Value = 479
$ws4.Range("F6").Value = 457
